# Update column G ("K" = strike count) values for the specified rows,
# per commit "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G8").Value = 1
